$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old "blank separator" row (row 11), pushing
# the separator row and the three summary rows down by two.
$ws.Rows.Item(11).Resize(2).Insert()

# New data row 11 (second shift on 2014-02-21, 09:00-10:45)
$ws.Range("A11").Value = 2014
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 21
$ws.Range("D11").Value = 0.375
$ws.Range("E11").Value = 0.44791666666666669
$ws.Range("D11").NumberFormat = "hh:mm;@"
$ws.Range("E11").NumberFormat = "hh:mm;@"
$ws.Range("G11").NumberFormat = "hh:mm;@"
$ws.Range("F11").NumberFormat = "0"

# New data row 12 (third shift on 2014-02-21, 12:15-13:00)
$ws.Range("A12").Value = 2014
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 21
$ws.Range("D12").Value = 0.51041666666666663
$ws.Range("E12").Value = 0.54166666666666663
$ws.Range("D12").NumberFormat = "hh:mm;@"
$ws.Range("E12").NumberFormat = "hh:mm;@"
$ws.Range("G12").NumberFormat = "hh:mm;@"
$ws.Range("F12").NumberFormat = "0"

# Extend the "time spent" shared formula down across the two new rows.
$ws.Range("F11:F12").Formula = "=(E11-D11)*24*60"

# Keep the blank separator row (now row 13) formatted like before.
$ws.Range("D13:G13").NumberFormat = "hh:mm;@"
$ws.Range("F13").NumberFormat = "0"

# Fix up the summary rows, now shifted to 14/15/16, to reference the new ranges.
$ws.Range("F14").Formula = "=SUM(F2:F13)"
$ws.Range("F15").Formula = "=F14/60"
$ws.Range("F16").Formula = "=F15/38.5"

[void]$ws.Range("F12").Select()
